# Auto-generated script applying scheduled market-data refresh
# to the per-job Leve profit sheets (columns H-N).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1455.0944
$ws.Range("J17").Value = 1176.3529
$ws.Range("L17").Value = 3529.0587
$ws.Range("N17").Value = -3865.0587
$ws.Range("H34").Value = 9999
$ws.Range("I34").Value = 9999
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 9999
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -9796
$ws.Range("N34").ClearContents()
$ws.Range("H36").Value = 9999
$ws.Range("I36").Value = 9999
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 9999
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -9284
$ws.Range("N36").ClearContents()
$ws.Range("H116").Value = 9763.385
$ws.Range("J116").Value = 2292.7
$ws.Range("L116").Value = 2292.7
$ws.Range("N116").Value = -9176.700000000001
$ws.Range("H140").Value = 56293.555
$ws.Range("J140").Value = 56293.555
$ws.Range("L140").Value = 56293.555
$ws.Range("N140").Value = -66653.55499999999
$ws.Range("H141").Value = 1575.8
$ws.Range("I141").Value = 802.25
$ws.Range("K141").Value = 2406.75
$ws.Range("M141").Value = 2773.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1852122
$ws.Range("I2").Value = 2778127.5
$ws.Range("K2").Value = 2778127.5
$ws.Range("M2").Value = -2778014.5
$ws.Range("H32").Value = 3141.2593
$ws.Range("I32").Value = 2898.9778
$ws.Range("K32").Value = 2898.9778
$ws.Range("M32").Value = -2611.9778
$ws.Range("H74").Value = 1737.3572
$ws.Range("I74").Value = 1592.4445
$ws.Range("K74").Value = 1592.4445
$ws.Range("M74").Value = -718.4445000000001
$ws.Range("H77").Value = 1737.3572
$ws.Range("I77").Value = 1592.4445
$ws.Range("K77").Value = 7962.2225
$ws.Range("M77").Value = -3594.2225
$ws.Range("H110").Value = 2598.6365
$ws.Range("I110").Value = 1861.4
$ws.Range("K110").Value = 1861.4
$ws.Range("M110").Value = 183.5999999999999
$ws.Range("H116").Value = 1852122
$ws.Range("I116").Value = 2778127.5
$ws.Range("K116").Value = 2778127.5
$ws.Range("M116").Value = -2775833.5
$ws.Range("H123").Value = 64998
$ws.Range("J123").Value = 64998
$ws.Range("L123").Value = 64998
$ws.Range("N123").Value = -74798
$ws.Range("H132").Value = 2479.08
$ws.Range("I132").Value = 1479.125
$ws.Range("K132").Value = 4437.375
$ws.Range("M132").Value = -1907.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1852122
$ws.Range("I3").Value = 2778127.5
$ws.Range("K3").Value = 2778127.5
$ws.Range("M3").Value = -2778013.5
$ws.Range("H20").Value = 2311.6296
$ws.Range("I20").Value = 2277.4092
$ws.Range("K20").Value = 2277.4092
$ws.Range("M20").Value = -2030.4092
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H86").Value = 185863.19
$ws.Range("I86").Value = 7872.25
$ws.Range("K86").Value = 7872.25
$ws.Range("M86").Value = -6749.25
$ws.Range("H89").Value = 185863.19
$ws.Range("I89").Value = 7872.25
$ws.Range("K89").Value = 39361.25
$ws.Range("M89").Value = -33745.25
$ws.Range("H92").Value = 19000
$ws.Range("J92").Value = 19000
$ws.Range("L92").Value = 19000
$ws.Range("N92").Value = -23992
$ws.Range("H107").Value = 1839.2273
$ws.Range("I107").Value = 1636.1818
$ws.Range("K107").Value = 1636.1818
$ws.Range("M107").Value = 283.8181999999999
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H134").Value = 8285.348
$ws.Range("I134").Value = 9008.684999999999
$ws.Range("K134").Value = 27026.055
$ws.Range("M134").Value = -24491.055

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1949.625
$ws.Range("I16").Value = 1119.6
$ws.Range("K16").Value = 1119.6
$ws.Range("M16").Value = -832.5999999999999
$ws.Range("H22").Value = 1620
$ws.Range("I22").Value = 100
$ws.Range("K22").Value = 100
$ws.Range("M22").Value = 250
$ws.Range("H31").Value = 1524.279
$ws.Range("I31").Value = 1533.75
$ws.Range("K31").Value = 1533.75
$ws.Range("M31").Value = -1238.75
$ws.Range("H34").Value = 1524.279
$ws.Range("I34").Value = 1533.75
$ws.Range("K34").Value = 1533.75
$ws.Range("M34").Value = -1331.75
$ws.Range("H58").Value = 3346175.2
$ws.Range("I58").Value = 3346175.2
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 3346175.2
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -3345972.2
$ws.Range("N58").ClearContents()
$ws.Range("H105").Value = 1302.6818
$ws.Range("I105").Value = 1050.375
$ws.Range("K105").Value = 1050.375
$ws.Range("M105").Value = 696.625
$ws.Range("H107").Value = 1086.421
$ws.Range("I107").Value = 1183.375
$ws.Range("J107").Value = 569.3333
$ws.Range("K107").Value = 1183.375
$ws.Range("L107").Value = 569.3333
$ws.Range("M107").Value = 736.625
$ws.Range("N107").Value = -4409.3333
$ws.Range("H113").Value = 1949.625
$ws.Range("I113").Value = 1119.6
$ws.Range("K113").Value = 1119.6
$ws.Range("M113").Value = 1050.4
$ws.Range("H122").Value = 4270.3335
$ws.Range("I122").Value = 1934.3334
$ws.Range("K122").Value = 5803.0002
$ws.Range("M122").Value = -3353.0002
$ws.Range("H136").Value = 3346175.2
$ws.Range("I136").Value = 3346175.2
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10038525.6
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -10035975.6
$ws.Range("N136").ClearContents()
$ws.Range("H141").Value = 71998.75
$ws.Range("J141").Value = 71998.75
$ws.Range("L141").Value = 71998.75
$ws.Range("N141").Value = -82358.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 672.6875
$ws.Range("I5").Value = 487.1111
$ws.Range("K5").Value = 1461.3333
$ws.Range("M5").Value = -1349.3333
$ws.Range("H26").Value = 607.4
$ws.Range("I26").Value = 782
$ws.Range("K26").Value = 2346
$ws.Range("M26").Value = -2058
$ws.Range("H98").Value = 849.3333
$ws.Range("I98").Value = 298
$ws.Range("J98").Value = 1125
$ws.Range("K98").Value = 894
$ws.Range("L98").Value = 3375
$ws.Range("M98").Value = 604
$ws.Range("N98").Value = -6371
$ws.Range("H131").Value = 12721.254
$ws.Range("J131").Value = 14072.528
$ws.Range("L131").Value = 42217.584
$ws.Range("N131").Value = -52297.584
$ws.Range("H135").Value = 672.6875
$ws.Range("I135").Value = 487.1111
$ws.Range("K135").Value = 4383.9999
$ws.Range("M135").Value = -1848.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1870.5385
$ws.Range("I97").Value = 1665.4546
$ws.Range("J97").Value = 2998.5
$ws.Range("K97").Value = 1665.4546
$ws.Range("L97").Value = 2998.5
$ws.Range("M97").Value = -1169.4546
$ws.Range("N97").Value = -3990.5
$ws.Range("H107").Value = 1470
$ws.Range("I107").Value = 76
$ws.Range("J107").Value = 2399.3333
$ws.Range("K107").Value = 76
$ws.Range("L107").Value = 2399.3333
$ws.Range("M107").Value = 1844
$ws.Range("N107").Value = -6239.3333
$ws.Range("H122").Value = 1348.1111
$ws.Range("I122").Value = 1019.1429
$ws.Range("K122").Value = 3057.4287
$ws.Range("M122").Value = -607.4287000000004
$ws.Range("H132").Value = 1101106.9
$ws.Range("I132").Value = 1426040.4
$ws.Range("K132").Value = 4278121.199999999
$ws.Range("M132").Value = -4275591.199999999
$ws.Range("H141").Value = 30331.666
$ws.Range("J141").Value = 30331.666
$ws.Range("L141").Value = 30331.666
$ws.Range("N141").Value = -40691.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1059.2858
$ws.Range("I82").Value = 827.5833
$ws.Range("J82").Value = 2449.5
$ws.Range("K82").Value = 827.5833
$ws.Range("L82").Value = 2449.5
$ws.Range("M82").Value = -466.5833
$ws.Range("N82").Value = -3171.5
$ws.Range("H85").Value = 1059.2858
$ws.Range("I85").Value = 827.5833
$ws.Range("J85").Value = 2449.5
$ws.Range("K85").Value = 827.5833
$ws.Range("L85").Value = 2449.5
$ws.Range("M85").Value = 420.4167
$ws.Range("N85").Value = -4945.5
$ws.Range("H132").Value = 1503.6666
$ws.Range("I132").Value = 953.6070999999999
$ws.Range("K132").Value = 2860.8213
$ws.Range("M132").Value = -330.8212999999996
$ws.Range("H136").Value = 2662.2122
$ws.Range("I136").Value = 1447.08
$ws.Range("K136").Value = 4341.24
$ws.Range("M136").Value = -1791.24

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 906.44446
$ws.Range("J107").Value = 1901.5
$ws.Range("L107").Value = 5704.5
$ws.Range("N107").Value = -9544.5
$ws.Range("H123").Value = 47599.918
$ws.Range("J123").Value = 47599.918
$ws.Range("L123").Value = 47599.918
$ws.Range("N123").Value = -57399.918
$ws.Range("H132").Value = 2076.4
$ws.Range("I132").Value = 1168.35
$ws.Range("K132").Value = 3505.05
$ws.Range("M132").Value = -975.0499999999997
